$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet name (title) to reflect the new "through" date
$ws.Name = "Through 2021-10-22"

# Update the label for the October row
$ws.Range("A11").Value = "October (through 10-22)"

# Update October row (row 11) values for columns C:H (2016-2021)
$ws.Range("C11").Value = 37
$ws.Range("D11").Value = 47
$ws.Range("E11").Value = 51
$ws.Range("F11").Value = 36
$ws.Range("G11").Value = 103
$ws.Range("H11").Value = 135

# Update Total row (row 12) values for columns C:H (2016-2021)
$ws.Range("C12").Value = 466
$ws.Range("D12").Value = 674
$ws.Range("E12").Value = 599
$ws.Range("F12").Value = 458
$ws.Range("G12").Value = 1004
$ws.Range("H12").Value = 1382
